$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.120179
$ws.Range("H2").Value = 0.360537
$ws.Range("I2").Value = 0.04921086431616203
$ws.Range("J2").Value = 0.04921086431616202
$ws.Range("M2").Value = 0.668273
$ws.Range("N2").Value = 2.004819
$ws.Range("O2").Value = 0.01328414746766746
$ws.Range("P2").Value = 0.01328414746766746
$ws.Range("Q2").Value = 0.080312380867
$ws.Range("R2").Value = 0.722811427803
$ws.Range("S2").Value = 0.0006537243785872707
$ws.Range("T2").Value = 0.0006537243785872705

# Row 3
$ws.Range("G3").Value = 0.120179
$ws.Range("H3").Value = 0.360537
$ws.Range("I3").Value = 0.04921086431616203
$ws.Range("J3").Value = 0.04921086431616202
$ws.Range("O3").Value = 0.3831531055114357
$ws.Range("P3").Value = 0.3831531055114357
$ws.Range("Q3").Value = 2.316440570620333
$ws.Range("R3").Value = 20.847965135583
$ws.Range("S3").Value = 0.01885529548763938
$ws.Range("T3").Value = 0.01885529548763937

# Row 4
$ws.Range("G4").Value = 0.120179
$ws.Range("H4").Value = 0.360537
$ws.Range("I4").Value = 0.04921086431616203
$ws.Range("J4").Value = 0.04921086431616202
$ws.Range("M4").Value = 30.36285833333334
$ws.Range("N4").Value = 91.08857500000001
$ws.Range("O4").Value = 0.6035627470208969
$ws.Range("P4").Value = 0.6035627470208967
$ws.Range("Q4").Value = 3.648977951641667
$ws.Range("R4").Value = 32.840801564775
$ws.Range("S4").Value = 0.02970184444993538
$ws.Range("T4").Value = 0.02970184444993537

# Row 5
$ws.Range("G5").Value = 1.522503666666667
$ws.Range("H5").Value = 4.567511000000001
$ws.Range("I5").Value = 0.6234343883806033
$ws.Range("J5").Value = 0.6234343883806033
$ws.Range("M5").Value = 0.668273
$ws.Range("N5").Value = 2.004819
$ws.Range("O5").Value = 0.01328414746766746
$ws.Range("P5").Value = 0.01328414746766746
$ws.Range("Q5").Value = 1.017448092834333
$ws.Range("R5").Value = 9.157032835509
$ws.Range("S5").Value = 0.008281794351663003
$ws.Range("T5").Value = 0.008281794351663001

# Row 6
$ws.Range("G6").Value = 1.522503666666667
$ws.Range("H6").Value = 4.567511000000001
$ws.Range("I6").Value = 0.6234343883806033
$ws.Range("J6").Value = 0.6234343883806033
$ws.Range("O6").Value = 0.3831531055114357
$ws.Range("P6").Value = 0.3831531055114357
$ws.Range("Q6").Value = 29.34613586720545
$ws.Range("R6").Value = 264.1152228048491
$ws.Range("S6").Value = 0.2388708219906507
$ws.Range("T6").Value = 0.2388708219906507

# Row 7
$ws.Range("G7").Value = 1.522503666666667
$ws.Range("H7").Value = 4.567511000000001
$ws.Range("I7").Value = 0.6234343883806033
$ws.Range("J7").Value = 0.6234343883806033
$ws.Range("M7").Value = 30.36285833333334
$ws.Range("N7").Value = 91.08857500000001
$ws.Range("O7").Value = 0.6035627470208969
$ws.Range("P7").Value = 0.6035627470208967
$ws.Range("Q7").Value = 46.22756314298056
$ws.Range("R7").Value = 416.0480682868251
$ws.Range("S7").Value = 0.3762817720382896
$ws.Range("T7").Value = 0.3762817720382896

# Row 8
$ws.Range("I8").Value = 0.3273547473032347
$ws.Range("J8").Value = 0.3273547473032347
$ws.Range("M8").Value = 0.668273
$ws.Range("N8").Value = 2.004819
$ws.Range("O8").Value = 0.01328414746766746
$ws.Range("P8").Value = 0.01328414746766746
$ws.Range("Q8").Value = 0.5342446126353333
$ws.Range("R8").Value = 4.808201513717999
$ws.Range("S8").Value = 0.004348628737417187
$ws.Range("T8").Value = 0.004348628737417186

# Row 9
$ws.Range("I9").Value = 0.3273547473032347
$ws.Range("J9").Value = 0.3273547473032347
$ws.Range("O9").Value = 0.3831531055114357
$ws.Range("P9").Value = 0.3831531055114357
$ws.Range("S9").Value = 0.1254269880331457
$ws.Range("T9").Value = 0.1254269880331456

# Row 10
$ws.Range("I10").Value = 0.3273547473032347
$ws.Range("J10").Value = 0.3273547473032347
$ws.Range("M10").Value = 30.36285833333334
$ws.Range("N10").Value = 91.08857500000001
$ws.Range("O10").Value = 0.6035627470208969
$ws.Range("P10").Value = 0.6035627470208967
$ws.Range("S10").Value = 0.1975791305326719
$ws.Range("T10").Value = 0.1975791305326718
